# Automatic update of files.
#
# The underlying species-occurrence export was re-synced upstream: the
# per-row "Taxonsorteringsordning" sort key (column B) was refreshed for
# every data row, which reshuffled which species record lands on which
# sheet row (rows 2, 4, 5, 6 and 7 now hold the records that used to sit
# on a different row; row 3's species stays put). Row 5 additionally picks
# up the "Publik kommentar" note ("Lodyta") that used to be on row 7.
#
# Cells are written individually (rather than copying whole row ranges) so
# that cells whose value doesn't actually change -- e.g. the always-blank
# I/K/AT placeholder cells present on every data row -- are left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- species previously on row 5, with a refreshed sort key.
$ws.Range("A2").Value = 111117771
$ws.Range("B2").Value = 99413
$ws.Range("D2").Value = "LC"
$ws.Range("E2").Value = 221235
$ws.Range("F2").Value = "Vårärt"
$ws.Range("G2").Value = "Lathyrus vernus"
$ws.Range("H2").Value = "(L.) Bernh."
$ws.Range("P2").Value = "Haxängberget (Haxängberget), Jmt"
$ws.Range("Q2").Value = 502162.4857801876
$ws.Range("R2").Value = 6984991.493007575

# Row 3 keeps its species; only the sort key is refreshed.
$ws.Range("B3").Value = 78623

# Row 4 <- species previously on row 6, with a refreshed sort key.
$ws.Range("A4").Value = 111115919
$ws.Range("B4").Value = 106732
$ws.Range("D4").Value = "NT"
$ws.Range("E4").Value = 220204
$ws.Range("F4").Value = "Slåtterfibbla"
$ws.Range("G4").Value = "Hypochaeris maculata"
$ws.Range("H4").Value = "L."
$ws.Range("Q4").Value = 502221.6631376348
$ws.Range("R4").Value = 6984971.081527092

# Row 5 <- species previously on row 7, with a refreshed sort key; also
# inherits the "Lodyta" public comment that was on the old row 7.
$ws.Range("A5").Value = 111115983
$ws.Range("B5").Value = 78578
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("P5").Value = "Haxängänget (Haxängänget), Jmt"
$ws.Range("Q5").Value = 502198.9066096816
$ws.Range("R5").Value = 6984972.883986787
$ws.Range("AC5").Value = "Lodyta"

# Row 6 <- species previously on row 4, with a refreshed sort key.
$ws.Range("A6").Value = 111116011
$ws.Range("B6").Value = 97565
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 2082
$ws.Range("F6").Value = "Skogsrör"
$ws.Range("G6").Value = "Calamagrostis chalybaea"
$ws.Range("H6").Value = "(Laest.) Fr."
$ws.Range("Q6").Value = 502202.5351659534
$ws.Range("R6").Value = 6984991.065364953

# Row 7 <- species previously on row 2, with a refreshed sort key; loses
# the "Lodyta" public comment (it moved to row 5 above), along with the
# now-unused Enhet/Metod/Bestämningsmetod placeholder cells that travelled
# with the old row 7 data.
$ws.Range("A7").Value = 111119600
$ws.Range("B7").Value = 96346
$ws.Range("E7").Value = 620
$ws.Range("F7").Value = "Skogsfru"
$ws.Range("G7").Value = "Epipogium aphyllum"
$ws.Range("H7").Value = "Sw."
$ws.Range("Q7").Value = 502198.48677184
$ws.Range("R7").Value = 6984920.618724592
$ws.Range("J7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("AC7").Value = ""
$ws.Range("AF7").Value = ""
